$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").Value = "test02-2"
[void]$ws.Range("A2").Select()
